$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J25
$values = @(
    @(6, 6),
    @(12, 13),
    @(5, 6),
    @(6, 7),
    @(1, 3),
    @(1, 6),
    @(1, 5),
    @(1, 3),
    @(1, 6),
    @(1, 6),
    @(1, 3),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(1, 4),
    @(5, 7),
    @(5, 6),
    @(1, 3),
    @(6, 6),
    @(1, 1),
    @(2, 2)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
